# Applies cryptos list price/volume refresh per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '30.641.26'; ForceText = $false }
    @{ Cell = 'D3'; Value = '1.874.36'; ForceText = $false }
    @{ Cell = 'E4'; Value = '  -0.05%  '; ForceText = $false }
    @{ Cell = 'D5'; Value = '247.72'; ForceText = $true }
    @{ Cell = 'E5'; Value = '  +1.05%  '; ForceText = $false }
    @{ Cell = 'E6'; Value = '  -0.04%  '; ForceText = $false }
    @{ Cell = 'E7'; Value = '  +0.09%  '; ForceText = $false }
    @{ Cell = 'D8'; Value = '0.2909'; ForceText = $true }
    @{ Cell = 'E8'; Value = '  +1.12%  '; ForceText = $false }
    @{ Cell = 'D9'; Value = '0.06486'; ForceText = $true }
    @{ Cell = 'E9'; Value = '  +0.26%  '; ForceText = $false }
    @{ Cell = 'D10'; Value = '22.15'; ForceText = $true }
    @{ Cell = 'E10'; Value = '  +5.06%  '; ForceText = $false }
    @{ Cell = 'D11'; Value = '0.07708'; ForceText = $true }
    @{ Cell = 'E11'; Value = '  -0.89%  '; ForceText = $false }
    @{ Cell = 'D12'; Value = '0.7387'; ForceText = $true }
    @{ Cell = 'E12'; Value = '  +1.15%  '; ForceText = $false }
    @{ Cell = 'D13'; Value = '96.41'; ForceText = $true }
    @{ Cell = 'E13'; Value = '  +1.34%  '; ForceText = $false }
    @{ Cell = 'D14'; Value = '1.873.21'; ForceText = $false }
    @{ Cell = 'E14'; Value = '  +0.12%  '; ForceText = $false }
    @{ Cell = 'D15'; Value = '5.165'; ForceText = $true }
    @{ Cell = 'E15'; Value = '  +0.44%  '; ForceText = $false }
    @{ Cell = 'D16'; Value = '273.67'; ForceText = $true }
    @{ Cell = 'E16'; Value = '  -0.55%  '; ForceText = $false }
    @{ Cell = 'D17'; Value = '30.656.67'; ForceText = $false }
    @{ Cell = 'E18'; Value = '  -0.34%  '; ForceText = $false }
    @{ Cell = 'E19'; Value = '  -0.03%  '; ForceText = $false }
    @{ Cell = 'D20'; Value = '0.000007515'; ForceText = $true }
    @{ Cell = 'E20'; Value = '  -0.49%  '; ForceText = $false }
    @{ Cell = 'D21'; Value = '2.118.49'; ForceText = $false }
    @{ Cell = 'E21'; Value = '  +0.31%  '; ForceText = $false }
    @{ Cell = 'D22'; Value = '1.000'; ForceText = $true }
    @{ Cell = 'E22'; Value = '  -0.02%  '; ForceText = $false }
    @{ Cell = 'D23'; Value = '5.263'; ForceText = $true }
    @{ Cell = 'E23'; Value = '  +0.39%  '; ForceText = $false }
    @{ Cell = 'D24'; Value = '6.196'; ForceText = $true }
    @{ Cell = 'E24'; Value = '  +0.45%  '; ForceText = $false }
    @{ Cell = 'D25'; Value = '9.218'; ForceText = $true }
    @{ Cell = 'E25'; Value = '  -0.47%  '; ForceText = $false }
    @{ Cell = 'D26'; Value = '163.93'; ForceText = $true }
    @{ Cell = 'E26'; Value = '  -0.93%  '; ForceText = $false }
    @{ Cell = 'D27'; Value = '18.76'; ForceText = $true }
    @{ Cell = 'E27'; Value = '  -0.67%  '; ForceText = $false }
    @{ Cell = 'D28'; Value = '1.912'; ForceText = $true }
    @{ Cell = 'E28'; Value = '  -0.11%  '; ForceText = $false }
    @{ Cell = 'D30'; Value = '1.346'; ForceText = $true }
    @{ Cell = 'E30'; Value = '  -2.57%  '; ForceText = $false }
    @{ Cell = 'D31'; Value = '1.508'; ForceText = $true }
    @{ Cell = 'E31'; Value = '  -1.00%  '; ForceText = $false }
    @{ Cell = 'D32'; Value = '4.281'; ForceText = $true }
    @{ Cell = 'E32'; Value = '  -0.63%  '; ForceText = $false }
    @{ Cell = 'D33'; Value = '4.102'; ForceText = $true }
    @{ Cell = 'E33'; Value = '  +1.41%  '; ForceText = $false }
    @{ Cell = 'D34'; Value = '0.04801'; ForceText = $true }
    @{ Cell = 'E34'; Value = '  +0.59%  '; ForceText = $false }
    @{ Cell = 'D35'; Value = '1.121'; ForceText = $true }
    @{ Cell = 'E35'; Value = '  -0.08%  '; ForceText = $false }
    @{ Cell = 'E36'; Value = '  -0.21%  '; ForceText = $false }
    @{ Cell = 'E37'; Value = '  +0.01%  '; ForceText = $false }
    @{ Cell = 'D38'; Value = '0.01856'; ForceText = $true }
    @{ Cell = 'E38'; Value = '  +0.68%  '; ForceText = $false }
    @{ Cell = 'D39'; Value = '2.754'; ForceText = $true }
    @{ Cell = 'E39'; Value = '  +0.11%  '; ForceText = $false }
    @{ Cell = 'D40'; Value = '6.245'; ForceText = $true }
    @{ Cell = 'E40'; Value = '  -2.55%  '; ForceText = $false }
    @{ Cell = 'D41'; Value = '73.13'; ForceText = $true }
    @{ Cell = 'E41'; Value = '  +4.22%  '; ForceText = $false }
    @{ Cell = 'D42'; Value = '1.973'; ForceText = $true }
    @{ Cell = 'E42'; Value = '  +3.01%  '; ForceText = $false }
    @{ Cell = 'D43'; Value = '0.4181'; ForceText = $true }
    @{ Cell = 'E43'; Value = '  +1.41%  '; ForceText = $false }
    @{ Cell = 'E44'; Value = '  +0.00%  '; ForceText = $false }
    @{ Cell = 'D45'; Value = '0.8365'; ForceText = $true }
    @{ Cell = 'E45'; Value = '  -0.46%  '; ForceText = $false }
    @{ Cell = 'D46'; Value = '102.41'; ForceText = $true }
    @{ Cell = 'E46'; Value = '  +0.10%  '; ForceText = $false }
    @{ Cell = 'D47'; Value = '9.326'; ForceText = $true }
    @{ Cell = 'E47'; Value = '  -0.55%  '; ForceText = $false }
    @{ Cell = 'D48'; Value = '35.46'; ForceText = $true }
    @{ Cell = 'E48'; Value = '  +0.56%  '; ForceText = $false }
    @{ Cell = 'D49'; Value = '6.988'; ForceText = $true }
    @{ Cell = 'E49'; Value = '  -1.32%  '; ForceText = $false }
    @{ Cell = 'D50'; Value = '918.93'; ForceText = $true }
    @{ Cell = 'E50'; Value = '  +0.15%  '; ForceText = $false }
    @{ Cell = 'D51'; Value = '0.05653'; ForceText = $true }
    @{ Cell = 'E51'; Value = '  +1.49%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $u.Value
}
